$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position tweak (xWindow 1185 -> 240) ---
$excel.ActiveWindow.Left = 240

# --- Append the two new news rows (row 4 and row 5) ---
$ws.Range("A4").Value = "https://tribune.com.pk/"
$ws.Range("B4").Value = "February 18, 2019"
$ws.Range("C4").Value = "Our Correspondent"
$ws.Range("D4").Value = "Dialogue only path to peace, stresses crown prince"
$ws.Range("E4").Value = "https://tribune.com.pk/story/1913278/1-believe-pakistan-crown-prince-mohammad-bin-salman-departs-pakistan/"

$ws.Range("A5").Value = "https://www.dawn.com"
$ws.Range("B5").Value = "February 17, 2019"
$ws.Range("C5").Value = "Anwar Iqbal"
$ws.Range("D5").Value = "Pulwama explosives obtained locally, says Indian commander`n"
$ws.Range("E5").Value = "https://www.dawn.com/news/1464358/pulwama-explosives-obtained-locally-says-indian-commander"

# Row 5's description contains an embedded line break; AutoFit keeps the row
# at its natural (default) height instead of leaving a stale explicit
# "wrapped height" behind, matching a plain unstyled <row> in the target.
$ws.Rows(5).AutoFit()

# --- Column widths: columns A and D get their own explicit width,
#     splitting them out of the wider merged column definitions
#     (A previously shared a width definition with B, D shared one with E) ---
$ws.Columns("A").ColumnWidth = 29.17
$ws.Columns("D").ColumnWidth = 49.17

# --- Reset the active selection back to A1 ---
$ws.Range("A1").Select() | Out-Null
